$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "epsilon algorithm": nudge the end time recorded for the 2014-02-22
# evening entry (row 20) down by 15 minutes (0.875 -> 0.85416666666666663).
# All dependent formulas (F20, G20, F22, F23, F24) recalc automatically.
$ws.Range("E20").Value = 0.85416666666666663

# Update the window/view state to match where the user ended up looking:
# scrolled so row 7 is the first visible row, with E21 as the active cell.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E21").Select()
